# Atualização de bases das ligas, do dia: 01-06-2024 às 01:16
#
# The source feed re-synced a handful of match rows and, in doing so,
# several row *pairs* (and one row *triple*) ended up with their entire
# record (every column except the running index in column A) exchanged
# with a sibling row. Column A keeps strictly increasing (0,1,2,...) as
# the row's positional index, so it is left untouched; everything from
# column B ("id") through column AD ("PL_AhUnder") moves as a unit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $rowA, $rowB) {
    $rangeA = $sheet.Range("B$rowA`:AD$rowA")
    $rangeB = $sheet.Range("B$rowB`:AD$rowB")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# Simple pairwise swaps: each pair exchanges its whole record.
Swap-Rows $ws 133 134
Swap-Rows $ws 229 231
Swap-Rows $ws 232 233
Swap-Rows $ws 307 308
Swap-Rows $ws 310 311

# Three-row rotation: row 263 receives row 264's record, row 264
# receives row 265's record, and row 265 receives row 263's (original)
# record.
$range263 = $ws.Range("B263:AD263")
$range264 = $ws.Range("B264:AD264")
$range265 = $ws.Range("B265:AD265")

$val263 = $range263.Value2
$val264 = $range264.Value2
$val265 = $range265.Value2

$range263.Value2 = $val264
$range264.Value2 = $val265
$range265.Value2 = $val263
